$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new date
$ws.Name = "Constraints_2018-07-31"

# Rename the constraint header from f_upper_lower_relation_bounds to
# height_end_start_relation_bounds
$ws.Range("F1").Value = "height_end_start_relation_bounds"

# Update the f_upper_lower_relation_bounds (now height_end_start_relation_bounds)
# and f_regression_percentage_bounds values for each pattern row

# Triangle
$ws.Range("F2").Value = "[]"
$ws.Range("G2").Value = "[]"

# Triangle top
$ws.Range("F3").Value = "[0.1, 0.5]"

# Triangle bottom
$ws.Range("F4").Value = "[0.1, 0.5]"

# Triangle up
$ws.Range("F5").Value = "[0.1, 0.5]"
$ws.Range("G5").Value = "[]"

# Triangle down
$ws.Range("F6").Value = "[0.1, 0.5]"
$ws.Range("G6").Value = "[]"

# Channel
$ws.Range("E7").Value = "[]"

# Channel up
$ws.Range("E8").Value = "[]"
$ws.Range("F8").Value = "[0.9, 1.1]"

# Channel down
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "[0.9, 1.1]"

# TKE up
$ws.Range("F10").Value = "[0.1, 0.5]"

# TKE down
$ws.Range("F11").Value = "[0.1, 0.5]"

# Head-Shoulder
$ws.Range("D12").Value = "[-1.0, 1.0]"
$ws.Range("E12").Value = "[]"
$ws.Range("F12").Value = "[]"
$ws.Range("G12").Value = "[-1.0, 1.0]"

# Inverse-Head-Shoulder
$ws.Range("D13").Value = "[-1.0, 1.0]"
$ws.Range("E13").Value = "[-1.0, 1.0]"
$ws.Range("F13").Value = "[]"
$ws.Range("G13").Value = "[-1.0, 1.0]"
